# Applies the "mostly ig generated files" regeneration diff to the
# StructureDefinition-inadvertent-administration-status workbook:
#   - Metadata!B2  URL: pythia -> cicada IG path
#   - Metadata!B8  Date: regenerated timestamp
#   - Metadata: new "Jurisdiction" row inserted after "Contact" (row 11),
#     pushing Description/Purpose/Copyright/... down by one row
#   - Elements!R5  (the Fixed Value column on the Extension.url row) carries
#     the same URL text as Metadata!B2 and must be kept in sync

$wb = $excel.ActiveWorkbook

$metadata = $wb.Worksheets.Item("Metadata")
$elements = $wb.Worksheets.Item("Elements")

# --- URL + Date updates -----------------------------------------------
$metadata.Range("B2").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/inadvertent-administration-status"
$metadata.Range("B8").Value = "2026-02-11T14:37:07-05:00"

$elements.Range("R5").Value = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/inadvertent-administration-status"

# --- Insert the new "Jurisdiction" property row ------------------------
# Row 11 was "Description"; a new blank-valued "Jurisdiction" row is
# inserted above it (matching the existing formatting of the data rows),
# shifting Description/Purpose/Copyright/FHIR Version/... down by one.
$metadata.Rows.Item(11).Insert()

$formatSource = $metadata.Range("A12:B12")
$newRow = $metadata.Range("A11:B11")
$formatSource.Copy()
$newRow.PasteSpecial(-4122)

$metadata.Range("A11").Value = "Jurisdiction"
$metadata.Range("B11").Value = ""
